$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Linear)
$ws.Range("B2").Value = 3.335376396475537
$ws.Range("C2").Value = 3.179886475244781
$ws.Range("D2").Value = 0.6645808329799155

# Row 3 (Decision Tree)
$ws.Range("B3").Value = 4.559463922686821
$ws.Range("C3").Value = 4.090476190476191
$ws.Range("D3").Value = 0.3732046855593157

# Row 4 (Random Forest)
$ws.Range("B4").Value = 4.191087595728871
$ws.Range("C4").Value = 3.518779342723005
$ws.Range("D4").Value = 0.4703955205505806

# Row 5 (Lasso)
$ws.Range("B5").Value = 1.888439805963124
$ws.Range("C5").Value = 1.55714761061008
$ws.Range("D5").Value = 0.8924762341483605

# Row 6 (Optimized Equation)
$ws.Range("B6").Value = 2.596390190944105
$ws.Range("C6").Value = 2.426800380288052
$ws.Range("D6").Value = 0.7967464716493238
